$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain-looking number (e.g. "315.42") must be
# force-written as TEXT (matching the source inlineStr cells) rather than
# letting Excel auto-convert them to the Number type. We do this by briefly
# switching the cell to the Text number format, assigning the value, then
# restoring the "Normal" style so no visible formatting/style change remains.

$ws.Range('D2').Value = '28.024.04'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '1.907.80'
$ws.Range('E3').Value = '  +2.28%  '
$ws.Range('E4').Value = '  -0.86%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('E7').Value = '  +0.78%  '
$ws.Range('E8').Value = '  +1.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07367'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.81'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07749'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.936.66'
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.497'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.649'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.005'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008831'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').Value = '28.056.06'
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.167'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('D23').Value = '2.141.60'
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('E24').Value = '  +1.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.922'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('E28').Value = '  +5.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '116.89'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.960'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08939'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.292'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.259'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7746'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.680'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.643'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02061'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.112'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05307'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5490'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.988'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.039'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1528'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.492'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4821'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '108.03'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.89%  '
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.29'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06075'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.08%  '
